# Restore functionality of unique_molecule counting function:
# update the heatmap values on Sheet1 to the corrected unique counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (pyridine)
$ws.Range("C2").Value = 2882
$ws.Range("D2").Value = 2842
$ws.Range("E2").Value = 1487
$ws.Range("F2").Value = 2882
$ws.Range("G2").Value = 2431
$ws.Range("H2").Value = 2251
$ws.Range("I2").Value = 2251
$ws.Range("J2").Value = 2251
$ws.Range("K2").Value = 2430
$ws.Range("L2").Value = 2882

# Row 3 (pyridazine)
$ws.Range("G3").Value = 190
$ws.Range("H3").Value = 190
$ws.Range("I3").Value = 190
$ws.Range("J3").Value = 190
$ws.Range("K3").Value = 190

# Row 4 (pyrimidine)
$ws.Range("H4").Value = 628
$ws.Range("I4").Value = 628
$ws.Range("J4").Value = 628

# Row 5 (pyrazine)
$ws.Range("G5").Value = 465
$ws.Range("H5").Value = 465
$ws.Range("I5").Value = 465
$ws.Range("J5").Value = 465
$ws.Range("K5").Value = 465
$ws.Range("L5").Value = 395

# Row 6 (pyrrole)
$ws.Range("B6").Value = 1038
$ws.Range("C6").Value = 520
$ws.Range("D6").Value = 1038
$ws.Range("E6").Value = 520
$ws.Range("G6").Value = 767
$ws.Range("H6").Value = 800
$ws.Range("I6").Value = 751
$ws.Range("J6").Value = 751
$ws.Range("K6").Value = 721
$ws.Range("L6").Value = 520

# Row 7 (pyrazole)
$ws.Range("B7").Value = 1365
$ws.Range("C7").Value = 1363
$ws.Range("D7").Value = 1365
$ws.Range("E7").Value = 1360
$ws.Range("I7").Value = 1365
$ws.Range("J7").Value = 1365
$ws.Range("K7").Value = 1365

# Row 8 (imidazole)
$ws.Range("B8").Value = 511
$ws.Range("C8").Value = 511
$ws.Range("D8").Value = 511
$ws.Range("E8").Value = 510
$ws.Range("I8").Value = 511
$ws.Range("J8").Value = 511
$ws.Range("K8").Value = 511
$ws.Range("L8").Value = 510

# Row 12 (furan)
$ws.Range("G12").Value = 767
$ws.Range("H12").Value = 861
$ws.Range("I12").Value = 861
$ws.Range("J12").Value = 861
$ws.Range("K12").Value = 767
